$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix diacritics ("afrikati") in several cells of the scenario description
$ws.Range("B8").Value  = "Neuspješna identifikacija pacijenta, hitni slucajevi"
$ws.Range("A18").Value = "6.Bira jedan od ponuđenih termina ili odustaje od zahtjeva"
$ws.Range("B19").Value = "7.Označava odabrani termin, ako je odabran, kao zauzet"
$ws.Range("B23").Value = "Verifikacija pacijenta neuspješna"
$ws.Range("B32").Value = "Pacijent je hitan slučaj"
$ws.Range("B35").Value = "1.Prima pacijenta u prvi mogući termin"
$ws.Range("A36").Value = "6.Bira ponuđeni termin ili odustaje od zahtjeva"
$ws.Range("B37").Value = "7.Označava odabrani termin, ako je odabran, kao zauzet"

# Update the view state to match the exported/scrolled position recorded after the edit
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B38").Select()
